$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header label for the new "table description" column. The cell
# inherits the row's existing formatting (style index 2), matching A1.
$ws.Range("B1").Value = "table description"

$ws.Range("B6").Select()
